$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHashes = @{
    2  = '$2b$10$ZEaFnEVEQI38wULA.FIu0.tKDH0jKtVtzuYpnhB178wasn.wIB35.'
    3  = '$2b$10$lmMyT4ygt.yM9AwW6svphOib69jMRlrFbzmLcHr/DpM1NBSBfdu2W'
    4  = '$2b$10$JcBFILFDaQ6/8l7fNMe.j.t0TUNP8o2.SZUlNDMnRIAmRVV6oLWda'
    5  = '$2b$10$EgE8ymvINud3RmuzK1gqZ.Dzn.VEVgJXTa23bqwKfA64bcZGllnfi'
    6  = '$2b$10$7Lty84joZ8e9arEUQJZzLeJvyJql.b1P.PgikPyLI.rbrJ6yOPIxG'
    7  = '$2b$10$jMTy.aT1lp68nw/5k3Nch.4eM0duCoE68HkW.IvA6/njcS/E4djia'
    8  = '$2b$10$o2NeTYIOOF4AuJsGONOIrOWLqQnQX5nmWiiiFeTK8aAfimHuzfGmy'
    9  = '$2b$10$.Mw.eubcNkG9kwUg6F4Ttu0m1.vSv268YJP44TdRXoIic2IOrNjmG'
    10 = '$2b$10$TdeIp9uu2CfuPU.EU8z.RenHpTr3y6xgBhorjX/dvaEDE/q/zurB6'
    11 = '$2b$10$bf5yzmJe7TvW3C2dNm0SHuJtBjS/f8x7Ll3PsL.Qsv9taTynTM4j.'
    12 = '$2b$10$5i2P7Qvf9Xmq8jb8d5wMceTY3QBgH346HyMFGu5dIzXRLAr6o2FlS'
    13 = '$2b$10$YUGe5U.T0V8JeIbcXjB6Ye.USepVA.U6yxr6qKuOBgFKvT/rwjg66'
    14 = '$2b$10$T8kw/A8J4tv/Y9sYQlHZi.wPqngiOtfC4D2WtA9jvqUauuOaBC.4y'
    15 = '$2b$10$.OBbHlvnwMqwuIwi.gDZGeRgpDgDCu5Xd9BeFhrYS0dDvILXHNMJy'
    16 = '$2b$10$x/KT0f.ERbvz94sFnA1tZuX/HHFFvPUCjZrXDvLL8asl175Umb8MS'
    17 = '$2b$10$gkgY0Va7ioozVHwAkB3AjOGV3AsZTGLZXr13EyKd7JNdI8qr16/3u'
    18 = '$2b$10$vWUVCX6OIm1xvqYhEgSsf.CZpecmdSivvSUQQKxwDNlNc431PQku2'
    19 = '$2b$10$S5lAV2nohy6k5pV8jsnq5Ou08z.IXiRZPyk6DG.X6SAsVS2WcHBvO'
    20 = '$2b$10$5QMDkzJaq1OIbAUduJyUa.tfCM49y527yikoAXNJI5hJEUrZSIyKi'
    21 = '$2b$10$mn8VBWmeIxbrvJpPvW9jjexm1ZHT1T6Rae.uiIeA.2jiETBHHsFgC'
    22 = '$2b$10$CMa3HfVaanJXM/RZ0z68Ku.Zgl26os6idVi0aeN0/od/9iZAK1XmG'
    23 = '$2b$10$g.ZsLVVjqH1M0y8Dfv.1s.gFEOG1oidNCgyULohKMn0HFIK1OCING'
    24 = '$2b$10$XtwjjqFCTAXEHgyL25PqNuovXBIU/ktlXS2CBng4MCT.68nA3nPyy'
    25 = '$2b$10$OulLiczTAVI3M5J.pdBqUOJA6FBneDmDPthC59Bs1Nhy3xUjzKnwe'
    26 = '$2b$10$pDJCcSl9R620anoYhpGG0uWqdin.mUi3py3gFReaQsH/6OkXnRMQC'
}

foreach ($row in $newHashes.Keys) {
    $ws.Range("C$row").Value = $newHashes[$row]
}
